$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
  2  = @('1.02','1.03055029072378','1.039582631608319','1.030233735645875','1.048059835917407','1.034988157074656','1.035690915942963','1.04236740876327','1.033045339321515','1.050820736706897')
  3  = @('1.02','1.031525713853264','1.040334047593937','1.031063231416374','1.04894409824415','1.035156002880944','1.036307707967147','1.042929228036171','1.033683090194355','1.05151678109444')
  4  = @('1.02','1.032157419879803','1.040820539874346','1.031600807132731','1.049516877315329','1.035263397277794','1.036706761791015','1.043292374736066','1.034095948083048','1.051967133040546')
  5  = @('1.02','1.032423117985252','1.04102512589282','1.031827002915326','1.049757816376129','1.035308254932028','1.036874510650812','1.043444947683673','1.034269558112291','1.052156451330735')
  6  = @('1.02','1.032467737399309','1.041059480537181','1.031864993811216','1.04979827943329','1.035315769652412','1.03690267561095','1.043470559814552','1.034298710584372','1.05218823810251')
  7  = @('1.02','1.032160969645517','1.040823273310393','1.031603828791742','1.04952009619578','1.035263997811875','1.036709003313208','1.04329441379364','1.034098267694888','1.051969662760315')
  8  = @('1.02','1.030879826778748','1.039836518200812','1.030513893640531','1.048358550956346','1.035045132198321','1.035899373860225','1.042557357909368','1.033260830174246','1.051055974749293')
  9  = @('1.02','1.028626479261632','1.038099903751741','1.028599752738654','1.04631643534116','1.034650196283176','1.034472348170658','1.041255646555638','1.031786661570097','1.049445722806263')
  10 = @('1.02','1.027127115180575','1.036943704306533','1.027328081870708','1.044958252949471','1.034380707133342','1.03352081485359','1.04038593452166','1.030804953030839','1.048372147136781')
  11 = @('1.02','1.026478564324758','1.036443440318482','1.026778499453484','1.044370928184025','1.034262551757269','1.033108758598985','1.04000890074484','1.030380128702086','1.047907274216036')
  12 = @('1.02','1.026237767138234','1.036257678137218','1.026574520515101','1.044152887763824','1.034218444024337','1.03295569803995','1.039868788100282','1.030222370356616','1.047734599740375')
  13 = @('1.02','1.026289414257257','1.036297522103178','1.026618267409103','1.044199652805379','1.034227915210475','1.032988530239534','1.039898845686638','1.03025620821683','1.047771638981723')
  14 = @('1.02','1.026458657838854','1.036428083964837','1.026761635210504','1.04435290247666','1.034258910274455','1.033096106652428','1.039997320320567','1.030367087519789','1.047893000882761')
  15 = @('1.02','1.026562948129963','1.036508535181066','1.026849990142725','1.044447340354313','1.034277978276887','1.033162387458885','1.040057985141569','1.030435409290407','1.047967775993636')
  16 = @('1.02','1.027170171834309','1.036976913258948','1.027364578275709','1.044997248252352','1.034388517897906','1.033548160987129','1.040410947769486','1.030833152824195','1.048402999161249')
  17 = @('1.02','1.027551250799699','1.037270816606135','1.027687650360092','1.045342399998524','1.034457464669391','1.03379013751062','1.040632234077622','1.031082717650552','1.04867600194962')
  18 = @('1.02','1.027773593627207','1.037442281837406','1.027876194978369','1.045543796065421','1.034497538789297','1.033931274784418','1.040761263853301','1.031228309800391','1.048835239048663')
  19 = @('1.02','1.027849417942646','1.037500753180103','1.027940501106425','1.045612479580635','1.034511179035249','1.033979398320029','1.040805252384992','1.031277957207392','1.048889534611862')
  20 = @('1.02','1.027510357809829','1.037239279785612','1.027652977209221','1.045305360726083','1.03445008195176','1.033764176079189','1.04060849656908','1.03105593909203','1.048646711409312')
  21 = @('1.02','1.026408816959709','1.036389635135028','1.026719412527169','1.044307770995501','1.034249789057672','1.033064428197647','1.0399683237957','1.030334435207677','1.04785726281772')
  22 = @('1.02','1.025716833580774','1.035855766857162','1.02613337225249','1.043681230994537','1.03412258687833','1.032624443229522','1.039565443226819','1.029881030871608','1.047360905879432')
  23 = @('1.02','1.026083610021508','1.036138748135349','1.026443954817898','1.044013306380636','1.034190139381379','1.032857689722368','1.039779053385268','1.030121366603583','1.047624033626545')
  24 = @('1.02','1.027528835388089','1.03725352980504','1.027668644198147','1.045322096947687','1.03445341832091','1.033775906943559','1.040619222659702','1.031068039092149','1.048659946546282')
  25 = @('1.02','1.029208521893179','1.038548594141008','1.029093830084578','1.04684380840346','1.034753391558863','1.034841304845954','1.041592510355074','1.032167585182034','1.049862029875621')
}

$cols = @('B','C','D','E','F','I','J','K','L','M')

foreach ($row in $data.Keys) {
  $values = $data[$row]
  for ($i = 0; $i -lt $cols.Length; $i++) {
    $col = $cols[$i]
    $val = [double]$values[$i]
    $ws.Range("$col$row").Value = $val
  }
}
